$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cell_F2 = @'
Here are the detailed steps from the video based on the frames provided:

1. **Open Weather Website**: The operator opens a weather website displaying the current temperature and weather conditions.

2. **View Temperature**: The temperature is displayed prominently on the website, showing the current temperature for a specific city.

3. **Spreadsheet Opened**: The operator has a spreadsheet open on the right side, with column A ready to receive city names and column B for temperatures.

4. **Select City**: The operator selects a city from the spreadsheet (column A).

5. **Search for Temperature**: The operator enters the city name into the weather website's search bar and initiates the search.

6. **View Results**: The website displays the current temperature and weather conditions for the searched city.

7. **Record Temperature**: The operator notes the temperature displayed on the website.

8. **Fill Spreadsheet**: The operator enters the recorded temperature into the corresponding cell in column B of the spreadsheet.

9. **Repeat Process**: The operator returns to the spreadsheet, selects the next city from column A, and repeats the search and recording process.

10. **Continue Until Complete**: This process continues for each city listed in column A, with temperatures being recorded in column B until all cities have been processed.

11. **Final Review**: Once all temperatures are recorded, the operator may review the spreadsheet for accuracy.

12. **Save Spreadsheet**: The operator saves the spreadsheet with the recorded temperatures for future reference. 

This sequence captures the main actions taken in the video.
'@
$ws.Range("F2").Value = $cell_F2

$cell_F3 = @'
To get the weather temperature for "Yunnan" following the steps shown in the video, you can follow these instructions:

1. **Open a Web Browser**: Start by launching your preferred web browser.

2. **Search for Weather**: In the search bar, type "Yunnan weather" and press Enter. This should lead you to a weather website that provides current weather information.

3. **Locate the Temperature**: Once the weather page for Yunnan loads, look for the current temperature displayed prominently on the page. It is usually shown in large numbers at the top.

4. **Record the Temperature**: Note down the temperature. If there are multiple readings (like high and low), make sure to record the one you need.

5. **Open Your Spreadsheet**: Switch to your spreadsheet application (like Google Sheets or Excel).

6. **Input the Temperature**: Click on the cell where you want to enter the temperature and type in the value you noted from the weather website.

7. **Repeat if Necessary**: If you need to check the weather for other locations, repeat the steps above for each location.

By following these steps, you can efficiently gather and record the weather temperature for Yunnan or any other city.
'@
$ws.Range("F3").Value = $cell_F3

$cell_F4 = @'
Based on the frames provided, here is the JSON output detailing the clicks and actions taken in the Excel video:

```json
{
  "actions": [
    {
      "step": 1,
      "clicked_cell": "A1",
      "cell_value": "",
      "timestamp": "00:00:00"
    },
    {
      "step": 2,
      "clicked_cell": "A2",
      "cell_value": "",
      "timestamp": "00:00:01"
    },
    {
      "step": 3,
      "clicked_cell": "A3",
      "cell_value": "",
      "timestamp": "00:00:02"
    },
    {
      "step": 4,
      "clicked_cell": "A4",
      "cell_value": "",
      "timestamp": "00:00:03"
    },
    {
      "step": 5,
      "clicked_cell": "A5",
      "cell_value": "",
      "timestamp": "00:00:04"
    },
    {
      "step": 6,
      "clicked_cell": "A6",
      "cell_value": "",
      "timestamp": "00:00:05"
    },
    {
      "step": 7,
      "clicked_cell": "A7",
      "cell_value": "",
      "timestamp": "00:00:06"
    },
    {
      "step": 8,
      "clicked_cell": "A8",
      "cell_value": "",
      "timestamp": "00:00:07"
    },
    {
      "step": 9,
      "clicked_cell": "A9",
      "cell_value": "",
      "timestamp": "00:00:08"
    },
    {
      "step": 10,
      "clicked_cell": "A10",
      "cell_value": "",
      "timestamp": "00:00:09"
    }
  ]
}
```

This JSON structure captures each click on the cells in column A, along with the corresponding timestamps. Adjust the timestamps as necessary based on the actual video time.
'@
$ws.Range("F4").Value = $cell_F4

$cell_F5 = @'
I can't identify the content of the images or the video frames you've provided. However, if you need assistance with a specific question or topic related to the frames, feel free to ask!
'@
$ws.Range("F5").Value = $cell_F5

$cell_G5 = @'
1. 成都市天气_百度搜索&#10;2.上海天气_百度搜索&#10;3.北京天气_百度搜索&#10;4.广州天气_百度搜索&#10;2. work.xlsx
'@
$ws.Range("G5").Value = $cell_G5

$cell_F6 = @'
Here’s a JSON representation of the actions based on the frames provided:

```json
[
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    },
    {
        "action": "click",
        "window_title": "Baidu Weather",
        "position": {"x": 200, "y": 300}
    }
]
```

**Note:** The pixel positions and specific actions are placeholders. You would need to replace them with the actual coordinates and actions based on the specific interactions in the video frames.
'@
$ws.Range("F6").Value = $cell_F6

$cell_F7 = @'
It looks like you're working with a video that demonstrates how to gather weather data for various cities and input it into an Excel sheet. Based on the frames you've provided, here's a simplified table format for the city names and their corresponding temperatures:

| City Name | Temperature |
|-----------|-------------|
| City 1    | 25°C        |
| City 2    | 35°C        |
| City 3    | 32°C        |
| City 4    | 35°C        |

(Note: Replace "City 1", "City 2", etc., with the actual city names as per your data collection process.) 

If you need further assistance or specific details, feel free to ask!
'@
$ws.Range("F7").Value = $cell_F7

$cell_G7 = @'
| City Name | Temperature |&#10;| 上海 | 35° |&#10;| 北京 | 32° |&#10;| 广州 | 35° |
'@
$ws.Range("G7").Value = $cell_G7

$cell_G8 = @'
1.上海&#10;2.北京&#10;3.广州
'@
$ws.Range("G8").Value = $cell_G8

$cell_F9 = @'
Based on the frames provided, here is the list of cities in the order they were processed:

1. City 1
2. City 2
3. City 3
4. City 4
5. City 5
6. City 6
7. City 7
8. City 8
9. City 9
10. City 10
11. City 11
12. City 12
13. City 13
14. City 14
15. City 15
16. City 16
17. City 17
18. City 18
19. City 19
20. City 20

(Note: The specific names of the cities are not provided in the frames, so they are labeled generically as "City 1," "City 2," etc.)
'@
$ws.Range("F9").Value = $cell_F9

$cell_G9 = @'
1.上海&#10;2.北京&#10;3.广州
'@
$ws.Range("G9").Value = $cell_G9

$cell_F13 = @'
```json
{
  "button": {
    "count": 4,
    "names": ["Add...", "Edit...", "Delete...", "Entitlements..."]
  },
  "combo box": {
    "count": 1,
    "names": ["Access Group:"]
  },
  "dropdown": {
    "count": 1,
    "names": ["Filter"]
  },
  "text box": {
    "count": 2,
    "names": ["ID", "Display Name"]
  },
  "label": {
    "count": 5,
    "names": ["Application Pools", "Filter", "ID", "Display Name", "Farm"]
  }
}
```
'@
$ws.Range("F13").Value = $cell_F13

$cell_E15 = @'
Analyze the provided video and learn only the steps and processes that are present in the video. For each step in the video, identify the step action type, action value, action description, next step, and the pixel coordinates of that step. Present the results in JSON format where each item represents a step and contains the fields: &#10;"step",&#10;"action_type", "action_value", "action_description",  "pixel_coordinates" and "next_step". Note that the resolution of the video when I recorded it is 2560*1379.
'@
$ws.Range("E15").Value = $cell_E15

$cell_F15 = @'
Based on the frames provided, here is the JSON representation of the steps and processes observed in the video:

```json
[
    {
        "step": 1,
        "action_type": "click",
        "action_value": "search",
        "action_description": "Clicked on the search bar to enter a location.",
        "pixel_coordinates": [100, 200],
        "next_step": 2
    },
    {
        "step": 2,
        "action_type": "input",
        "action_value": "location name",
        "action_description": "Entered the name of the location in the search bar.",
        "pixel_coordinates": [150, 200],
        "next_step": 3
    },
    {
        "step": 3,
        "action_type": "submit",
        "action_value": "search",
        "action_description": "Submitted the search for the location.",
        "pixel_coordinates": [200, 200],
        "next_step": 4
    },
    {
        "step": 4,
        "action_type": "view",
        "action_value": "weather data",
        "action_description": "Viewed the weather data for the searched location.",
        "pixel_coordinates": [250, 200],
        "next_step": 5
    },
    {
        "step": 5,
        "action_type": "click",
        "action_value": "next day",
        "action_description": "Clicked to view the weather for the next day.",
        "pixel_coordinates": [300, 200],
        "next_step": 6
    },
    {
        "step": 6,
        "action_type": "view",
        "action_value": "next day weather",
        "action_description": "Viewed the weather data for the next day.",
        "pixel_coordinates": [350, 200],
        "next_step": 7
    },
    {
        "step": 7,
        "action_type": "click",
        "action_value": "details",
        "action_description": "Clicked to view detailed weather information.",
        "pixel_coordinates": [400, 200],
        "next_step": 8
    },
    {
        "step": 8,
        "action_type": "view",
        "action_value": "detailed weather",
        "action_description": "Viewed detailed weather information.",
        "pixel_coordinates": [450, 200],
        "next_step": 9
    },
    {
        "step": 9,
        "action_type": "click",
        "action_value": "back",
        "action_description": "Clicked to go back to the main weather page.",
        "pixel_coordinates": [500, 200],
        "next_step": 10
    },
    {
        "step": 10,
        "action_type": "click",
        "action_value": "another location",
        "action_description": "Clicked to search for another location.",
        "pixel_coordinates": [550, 200],
        "next_step": 11
    },
    {
        "step": 11,
        "action_type": "input",
        "action_value": "new location",
        "action_description": "Entered a new location in the search bar.",
        "pixel_coordinates": [600, 200],
        "next_step": 12
    },
    {
        "step": 12,
        "action_type": "submit",
        "action_value": "search",
        "action_description": "Submitted the search for the new location.",
        "pixel_coordinates": [650, 200],
        "next_step": 13
    },
    {
        "step": 13,
        "action_type": "view",
        "action_value": "new weather data",
        "action_description": "Viewed the weather data for the new location.",
        "pixel_coordinates": [700, 200],
        "next_step": 14
    }
]
```

This JSON structure outlines the steps taken in the video, including the action type, value, description, pixel coordinates, and the next step in the process. Adjust the pixel coordinates as necessary based on the actual positions in your frames.
'@
$ws.Range("F15").Value = $cell_F15

$cell_G17 = @'
```json&#10;{&#10;    "A1": "too",&#10;    "B1": "t00",&#10;    "C1": "ll",&#10;    "D1": "iill",&#10;    "E1": "LL11ll",&#10;    "A3": "top",&#10;    "B3": "100",&#10;    "C3": "lily",&#10;    "E4": "Lee"&#10;}&#10;```
'@
$ws.Range("G17").Value = $cell_G17

$cell_F18 = @'
Here’s the JSON representation of the tables in the image, focusing on the specified details:

```json
{
  "CurrentUsersSummary": {
    "attributes": ["User Id", "Home Branch", "Current Branch"],
    "values": [
      ["FCU50100", "001", "000"],
      ["ARUNR01", "001", "000"],
      ["INFRA32", "001", "001"],
      ["ASIN12", "001", "002"]
    ]
  },
  "UserAlerts": {
    "attributes": ["User ID", "Sequence Number", "Alert Type"],
    "values": [
      ["31618401", "31618401", "A"],
      ["31618401", "31618402", "D"],
      ["31618401", "31618403", "A"]
    ]
  },
  "UserRoleAccountClass": {
    "attributes": ["Role Id", "Account Class"],
    "values": [
      ["D", "D"]
    ]
  },
  "UserRoles": {
    "attributes": ["Role Description", "Branches Allowed", "Account Class"],
    "values": [
      ["FUNDS TRANSFER TxD", "D", "D"],
      ["FUNDS TRANSFER TxD", "D", "D"],
      ["Foreign Exchange Main D", "D", "D"]
    ]
  },
  "ModuleDashboardSummary": {
    "attributes": ["Module Id", "Module Description"],
    "values": [
      ["AC", "Accounting"],
      ["AD", "Auto of Day"],
      ["AM", "Asset Management"],
      ["BC", "Bills and Collections"],
      ["CM", "Card Maintenance"]
    ]
  }
}
```

This JSON structure captures the tables, their attributes, and values as seen in the image, while ensuring accuracy in recognizing characters.
'@
$ws.Range("F18").Value = $cell_F18

$cell_F19 = @'
To merge the first row in your Excel spreadsheet, follow these detailed steps based on the current state shown in the image:

1. **Select the Cells to Merge**:
   - Click and drag to select the cells in the first row that you want to merge. For example, if you want to merge cells A1, B1, and C1, click on cell A1 and drag to cell C1.

2. **Open the Merge Options**:
   - After selecting the cells, look for the "Merge & Center" option in the toolbar. This is usually found in the "Home" tab of the ribbon at the top of the Excel window.

3. **Choose the Merge Type**:
   - Click on the "Merge & Center" dropdown arrow to see the merge options. You can choose:
     - **Merge & Center**: Merges the selected cells and centers the content.
     - **Merge Across**: Merges cells in each row individually (not applicable if you are merging a single row).
     - **Merge Cells**: Merges the selected cells without centering the content.
     - **Unmerge Cells**: If you need to unmerge previously merged cells.

4. **Confirm the Merge**:
   - Click on your desired merge option (e.g., "Merge & Center"). The selected cells will now be merged into one cell.

5. **Adjust Formatting (if needed)**:
   - After merging, you may want to adjust the font size, style, or alignment of the text in the merged cell to ensure it looks the way you want.

By following these steps, you will successfully merge the first row in your Excel spreadsheet.
'@
$ws.Range("F19").Value = $cell_F19

$cell_D20 = @'
input/video/2.mp4
'@
$ws.Range("D20").Value = $cell_D20

$cell_F20 = @'
It looks like you've shared a series of frames from a video, possibly showing an Excel spreadsheet. If you need help with something specific regarding these frames, such as analyzing data or extracting information, please let me know!
'@
$ws.Range("F20").Value = $cell_F20

$cell_D21 = @'
input/video/3.mp4
'@
$ws.Range("D21").Value = $cell_D21

$cell_F21 = @'
To determine the number of click events based on the visible water ripples in the frames you provided, I would need to analyze each frame for the presence of these ripples. 
Please count the number of frames that show water ripples, as each ripple typically indicates a click event. If you can provide that count, I can help you summarize the findings or assist with any further analysis.
'@
$ws.Range("F21").Value = $cell_F21

$cell_D22 = @'
input/video/3.mp4
'@
$ws.Range("D22").Value = $cell_D22

$cell_F22 = @'
To determine the number of click events based on the visible selection frames in the provided images, we can analyze each frame for the presence of a selection frame around a cell.
From the frames you provided, it appears that there are **12 visible selection frames**. Each frame indicates a click event in Excel.
Therefore, the total number of clicks is **12**.
'@
$ws.Range("F22").Value = $cell_F22

$cell_D23 = @'
input/video/4.mp4
'@
$ws.Range("D23").Value = $cell_D23

$cell_F23 = @'
In the provided frames, the prompt box with a red background and white text containing the words "<左键>" appears **five times**. Each occurrence represents a left-click event.
'@
$ws.Range("F23").Value = $cell_F23
